$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for a new first column: shift the existing header row
# (gr number / username / password, currently in A1:C1) one cell to the
# right so it lands in B1:D1, leaving A1 free for the new "name" column.
$ws.Range("A1").Insert(-4161)  # xlShiftToRight

# Row 1: new header
$ws.Range("A1").Value = "name"

# Row 2: dummy data
$ws.Range("A2").Value = "m"
$ws.Range("B2").NumberFormat = "@"   # keep "1" as text, not a number
$ws.Range("B2").Value = "1"
$ws.Range("C2").Value = "m"
$ws.Range("D2").Value = "p"

# Row 3: dummy data
$ws.Range("A3").Value = "r"
$ws.Range("B3").NumberFormat = "@"   # keep "2" as text, not a number
$ws.Range("B3").Value = "2"
$ws.Range("C3").Value = "r"
$ws.Range("D3").Value = "r"
